$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update order number ("เลขที่ใบสั่งขาย") in column A from 2100000066 to
# 2100000073 for rows 2 through 9, while preserving the original cell
# formatting (the cells use a quote-prefixed numeric style). Assigning
# .Value directly would re-derive a fresh (unused) style without the
# quote-prefix flag, so we stash the existing format in a scratch cell,
# write the new value, then restore the captured format.
$scratch = $ws.Range("L1")

for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 1)

    $cell.Copy()
    $scratch.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $cell.Value = 2100000073

    $scratch.Copy()
    $cell.PasteSpecial(-4122) | Out-Null      # xlPasteFormats
}

$scratch.Clear()

# Update the active selection on the sheet to B13:C13 (active cell B13)
$ws.Range("B13:C13").Select()
